# TrendyOwls.xlsx — "Add files via upload" edit
#
# The refreshed data drops the "Dereck Lively II" row entirely and moves
# the "Tyrese Haliburton" row up so it sits right after "Mark Williams"
# (instead of after "Karl-Anthony Towns"). Net effect: the sheet shrinks
# from 17 data rows (+1 header) to 16 data rows (+1 header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row above the current row 14 ("OG Anunoby"), shifting
#    OG Anunoby / Karl-Anthony Towns / Tyrese Haliburton / Franz Wagner /
#    Dereck Lively II down by one row each.
$ws.Rows.Item(14).Insert()

# 2) Fill the newly-inserted row 14 with Tyrese Haliburton's record (this
#    is the row that used to live at row 16, now shifted to row 17).
$ws.Cells.Item(14, 1).Value = "Tyrese Haliburton"
$ws.Cells.Item(14, 2).Value = "PG,SG"
$ws.Cells.Item(14, 3).Value = "Indiana Pacers"

# 3) Remove the old Tyrese Haliburton row, now a duplicate sitting at row 17.
$ws.Rows.Item(17).Delete()

# 4) Remove the trailing "Dereck Lively II" row, now at row 18.
$ws.Rows.Item(18).Delete()
